$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block -----------------------------------------------------
# Invoice date moved forward one month (19-Mar-2024 -> 19-Apr-2024).
# Assign the raw serial number so the existing date NumberFormat/style
# on the cell is left untouched.
$ws.Range("J4").Value = 45401

# Bill-to block: swap out the old placeholder client for the real one.
# "Yee" (B5) -> "Ruwaida" now lives directly in A5; B5 is cleared.
$ws.Range("A5").Value = "Ruwaida"
$ws.Range("B5").Value = $null

# "ABN\ACN: 123456" (B6) -> company name now lives directly in A6; B6 cleared.
$ws.Range("A6").Value = "UR Refrigeration & Air Conditioning"
$ws.Range("B6").Value = $null

# New address line added under the company name.
$ws.Range("A7").Value = "24-28 Adderley Street East. Lidcombe. NSW. 2141"

# --- Job / site address -------------------------------------------------
$ws.Range("A12").Value = "70-72 Castle Hill Road, West Pennant Hills"

# --- Line items ----------------------------------------------------------
# Row 14 description is unchanged ("Mechanical Service design and documentation").

# Row 15 gains a brand-new line item with amounts.
$ws.Range("A15").Value = "Mechanical Service certification of car park makeup"
$ws.Range("G15").Value = 1000
$ws.Range("H15").Value = 1000
$ws.Range("I15").Value = 1000
$ws.Range("J15").Value = 1100

# Row 16 used to hold the "Fire Service design and documentation" line;
# it is now repurposed as the bold "Total:" row.
$ws.Range("A16").Value = $null
$ws.Range("F16").Value = "Total: "
$ws.Range("F16").Font.Bold = $true
$ws.Range("G16").Value = 1000
$ws.Range("G16").Font.Bold = $true

# --- Totals footer ---------------------------------------------------
$ws.Range("I39").Value = 1000
$ws.Range("J39").Value = 1100
